# NIT-9010916544.xlsx update: "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# - Removes 3 workers who are no longer in arrears (MARIA GRACIELA PALACIO
#   VEGA, ELKIN ALBERTO NUÑEZ SORACA, EDINSON MANUEL GOMEZ OLIVARES) by
#   deleting their entire table rows.
# - Rolls the arrears period forward from 2507 to 2508 for the remaining
#   workers.
# - Updates the Salario Basico for CARLOS ANDRES ACOSTA TERAN and DAIRIS
#   MILETH PINEDA ROJAS.
# - Refreshes the account summary totals (VALOR MORA, Cant. Trabajadores).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three workers that are no longer in the account statement.
# Deleting bottom-up keeps the remaining row numbers stable while iterating.
$ws.Rows("21").Delete()
$ws.Rows("19").Delete()
$ws.Rows("17").Delete()

# Roll the arrears period forward for the remaining workers (rows 16-20).
$ws.Range("E16:E20").Value = "2508"

# Updated Salario Basico values.
$ws.Range("G16").Value = 1579000
$ws.Range("G18").Value = 1423500

# Refreshed account summary.
$ws.Range("E11").Value = 290920
$ws.Range("C13").Value = 5
